$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay text (inline/shared strings),
# matching the source data which never uses real Number cells. Forcing the NumberFormat to
# "@" (Text) before assigning the value prevents Excel from silently re-interpreting values
# such as "216.48" as a floating point number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.165.33"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.86"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.72"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.20"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.644.55"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.47"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.162.29"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.48"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.53"
$ws.Range("E22").Value = "  +4.26%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.69"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.301.94"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.853"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.776.32"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.12"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.23"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.64"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -0.77%  "
